# Sample Project / Main.xlsx - "SAVE" edit
# The rule row for "R40" (row 11) has its Rule-name cell (column B)
# changed from the text "R40" to the text "1". The leading single quote
# forces the new value to be stored as text (matching the original
# shared-string / t="s" cell type) instead of being auto-detected as a
# number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
